$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying B:G data (ME, MAE, MSE, RMSE, SE, N) for each quarter label
# shifted down by one row (rows 2-10 -> rows 3-11), and a brand-new row of
# values was computed for row 2. Row labels in column A (Q0..Q9) stay fixed.

# Capture existing data for rows 2-10 (B:G) before overwriting.
$oldValues = @{}
for ($r = 2; $r -le 10; $r++) {
    $oldValues[$r] = $ws.Range("B$r" + ":G$r").Value2
}

# Shift old row r data down into row r+1 (process bottom-up to avoid clobbering).
for ($r = 10; $r -ge 2; $r--) {
    $ws.Range("B$($r+1)" + ":G$($r+1)").Value2 = $oldValues[$r]
}

# Write the newly computed values into row 2.
$ws.Range("B2").Value2 = -0.1338314788754218
$ws.Range("C2").Value2 = 2.04803740314168
$ws.Range("D2").Value2 = 20.30227337679758
$ws.Range("E2").Value2 = 4.505804409514197
$ws.Range("F2").Value2 = 4.605038435222934
$ws.Range("G2").Value2 = 23
